$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.836.59"
Set-TextValue "E2" "  +2.03%  "
Set-TextValue "D3" "1.890.77"
Set-TextValue "E3" "  +2.73%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "238.36"
Set-TextValue "E5" "  +2.50%  "
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "D7" "0.4801"
Set-TextValue "E7" "  +2.83%  "
Set-TextValue "D8" "0.2857"
Set-TextValue "E8" "  +5.28%  "
Set-TextValue "D9" "0.06527"
Set-TextValue "D10" "18.92"
Set-TextValue "E10" "  +17.60%  "
Set-TextValue "D11" "1.886.91"
Set-TextValue "E11" "  +2.54%  "
Set-TextValue "D12" "0.07577"
Set-TextValue "E12" "  +2.10%  "
Set-TextValue "D13" "94.83"
Set-TextValue "E13" "  +13.33%  "
Set-TextValue "D14" "5.135"
Set-TextValue "E14" "  +4.11%  "
Set-TextValue "D15" "0.6536"
Set-TextValue "E15" "  +5.51%  "
Set-TextValue "D16" "299.11"
Set-TextValue "E16" "  +32.46%  "
Set-TextValue "D17" "30.815.98"
Set-TextValue "E17" "  +2.21%  "
Set-TextValue "D18" "13.23"
Set-TextValue "E18" "  +7.22%  "
Set-TextValue "D19" "1.001"
Set-TextValue "E19" "  +0.09%  "
Set-TextValue "D20" "0.000007538"
Set-TextValue "E20" "  +3.66%  "
Set-TextValue "D21" "2.138.46"
Set-TextValue "E21" "  +2.76%  "
Set-TextValue "D22" "1.001"
Set-TextValue "E22" "  -0.14%  "
Set-TextValue "D23" "5.175"
Set-TextValue "E23" "  +5.79%  "
Set-TextValue "D24" "6.165"
Set-TextValue "E24" "  +5.38%  "
Set-TextValue "D25" "9.330"
Set-TextValue "E25" "  +1.61%  "
Set-TextValue "D26" "168.03"
Set-TextValue "E26" "  +2.28%  "
Set-TextValue "D27" "19.67"
Set-TextValue "E27" "  +10.80%  "
Set-TextValue "D28" "1.955"
Set-TextValue "E28" "  +4.99%  "
Set-TextValue "D29" "0.1068"
Set-TextValue "E29" "  +3.39%  "
Set-TextValue "D30" "1.363"
Set-TextValue "E30" "  -0.76%  "
Set-TextValue "D31" "4.210"
Set-TextValue "E31" "  +3.36%  "
Set-TextValue "D32" "3.973"
Set-TextValue "E32" "  +4.38%  "
Set-TextValue "D33" "0.05034"
Set-TextValue "E33" "  +4.47%  "
Set-TextValue "E34" "  +2.84%  "
Set-TextValue "D35" "0.7257"
Set-TextValue "D36" "2.718"
Set-TextValue "E36" "  +0.52%  "
Set-TextValue "D37" "0.01951"
Set-TextValue "E37" "  +4.40%  "
Set-TextValue "D38" "2.719"
Set-TextValue "E38" "  +2.61%  "
Set-TextValue "D39" "2.060"
Set-TextValue "E39" "  +7.39%  "
Set-TextValue "D40" "0.8981"
Set-TextValue "E40" "  +0.70%  "
Set-TextValue "D41" "107.70"
Set-TextValue "E41" "  +3.40%  "
Set-TextValue "D42" "1.002"
Set-TextValue "E42" "  -0.06%  "
Set-TextValue "D43" "0.4209"
Set-TextValue "E43" "  +5.02%  "
Set-TextValue "D44" "5.603"
Set-TextValue "E44" "  +1.50%  "
Set-TextValue "D45" "66.45"
Set-TextValue "E45" "  +11.27%  "
Set-TextValue "D46" "7.371"
Set-TextValue "E46" "  +4.85%  "
Set-TextValue "D47" "0.1230"
Set-TextValue "E47" "  +3.22%  "
Set-TextValue "D48" "8.905"
Set-TextValue "E48" "  +3.54%  "
Set-TextValue "D49" "34.67"
Set-TextValue "E49" "  +5.47%  "
Set-TextValue "D50" "0.05626"
Set-TextValue "E50" "  +2.04%  "
Set-TextValue "E51" "  +2.75%  "

Write-Output "Applied updates to $($ws.Name)"
